$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update F-column query timestamps on the "data" sheet ---
$newTimestamps = @(
    "2021-10-05 14:34:00.666979",
    "2021-10-05 14:34:00.666987",
    "2021-10-05 14:34:00.666991",
    "2021-10-05 14:34:00.666993",
    "2021-10-05 14:34:00.666997",
    "2021-10-05 14:34:00.667000",
    "2021-10-05 14:34:00.667003",
    "2021-10-05 14:34:00.667005",
    "2021-10-05 14:34:00.667008",
    "2021-10-05 14:34:00.667011",
    "2021-10-05 14:34:00.667014",
    "2021-10-05 14:34:00.667017",
    "2021-10-05 14:34:00.667019",
    "2021-10-05 14:34:00.667022",
    "2021-10-05 14:34:00.667025",
    "2021-10-05 14:34:00.667027",
    "2021-10-05 14:34:00.667030",
    "2021-10-05 14:34:00.667033",
    "2021-10-05 14:34:00.667036",
    "2021-10-05 14:34:00.667038",
    "2021-10-05 14:34:00.667041",
    "2021-10-05 14:34:00.667044",
    "2021-10-05 14:34:00.667046",
    "2021-10-05 14:34:00.667049",
    "2021-10-05 14:34:00.667052",
    "2021-10-05 14:34:00.667055",
    "2021-10-05 14:34:00.667057",
    "2021-10-05 14:34:00.667060",
    "2021-10-05 14:34:00.667063",
    "2021-10-05 14:34:00.667065",
    "2021-10-05 14:34:00.667068",
    "2021-10-05 14:34:00.667071",
    "2021-10-05 14:34:00.667074",
    "2021-10-05 14:34:00.667077",
    "2021-10-05 14:34:00.667080",
    "2021-10-05 14:34:00.667082",
    "2021-10-05 14:34:00.667085",
    "2021-10-05 14:34:00.667088",
    "2021-10-05 14:34:00.667090",
    "2021-10-05 14:34:00.667093",
    "2021-10-05 14:34:00.667096",
    "2021-10-05 14:34:00.667099",
    "2021-10-05 14:34:00.667102",
    "2021-10-05 14:34:00.667104",
    "2021-10-05 14:34:00.667107",
    "2021-10-05 14:34:00.667110",
    "2021-10-05 14:34:00.667113",
    "2021-10-05 14:34:00.667115",
    "2021-10-05 14:34:00.667118",
    "2021-10-05 14:34:00.667121",
    "2021-10-05 14:34:00.667124",
    "2021-10-05 14:34:00.667126",
    "2021-10-05 14:34:00.667129",
    "2021-10-05 14:34:00.667132",
    "2021-10-05 14:34:00.667135",
    "2021-10-05 14:34:00.667137",
    "2021-10-05 14:34:00.667141",
    "2021-10-05 14:34:00.667143",
    "2021-10-05 14:34:00.667146",
    "2021-10-05 14:34:00.667149",
    "2021-10-05 14:34:00.667151",
    "2021-10-05 14:34:00.667154",
    "2021-10-05 14:34:00.667157",
    "2021-10-05 14:34:00.667160",
    "2021-10-05 14:34:00.667164",
    "2021-10-05 14:34:00.667167",
    "2021-10-05 14:34:00.667170",
    "2021-10-05 14:34:00.667172",
    "2021-10-05 14:34:00.667175",
    "2021-10-05 14:34:00.667178",
    "2021-10-05 14:34:00.667181",
    "2021-10-05 14:34:00.667183",
    "2021-10-05 14:34:00.667186",
    "2021-10-05 14:34:00.667189",
    "2021-10-05 14:34:00.667191",
    "2021-10-05 14:34:00.667194",
    "2021-10-05 14:34:00.667199",
    "2021-10-05 14:34:00.667202",
    "2021-10-05 14:34:00.667205",
    "2021-10-05 14:34:00.667208",
    "2021-10-05 14:34:00.667210",
    "2021-10-05 14:34:00.667213",
    "2021-10-05 14:34:00.667216",
    "2021-10-05 14:34:00.667218",
    "2021-10-05 14:34:00.667221",
    "2021-10-05 14:34:00.667224",
    "2021-10-05 14:34:00.667227",
    "2021-10-05 14:34:00.667229",
    "2021-10-05 14:34:00.667232",
    "2021-10-05 14:34:00.667234",
    "2021-10-05 14:34:00.667237",
    "2021-10-05 14:34:00.667240",
    "2021-10-05 14:34:00.667243",
    "2021-10-05 14:34:00.667246",
    "2021-10-05 14:34:00.667249",
    "2021-10-05 14:34:00.667251",
    "2021-10-05 14:34:00.667254",
    "2021-10-05 14:34:00.667256",
    "2021-10-05 14:34:00.667259",
    "2021-10-05 14:34:00.667261",
    "2021-10-05 14:34:00.667264",
    "2021-10-05 14:34:00.667267",
    "2021-10-05 14:34:00.667269",
    "2021-10-05 14:34:00.667272",
    "2021-10-05 14:34:00.667274",
    "2021-10-05 14:34:00.667277",
    "2021-10-05 14:34:00.667280",
    "2021-10-05 14:34:00.667282",
    "2021-10-05 14:34:00.667287",
    "2021-10-05 14:34:00.667290",
    "2021-10-05 14:34:00.667293",
    "2021-10-05 14:34:00.667295",
    "2021-10-05 14:34:00.667298",
    "2021-10-05 14:34:00.667301",
    "2021-10-05 14:34:00.667303",
    "2021-10-05 14:34:00.667306",
    "2021-10-05 14:34:00.667309",
    "2021-10-05 14:34:00.667311",
    "2021-10-05 14:34:00.667314",
    "2021-10-05 14:34:00.667317",
    "2021-10-05 14:34:00.667319",
    "2021-10-05 14:34:00.667322",
    "2021-10-05 14:34:00.667325",
    "2021-10-05 14:34:00.667327",
    "2021-10-05 14:34:00.667330",
    "2021-10-05 14:34:00.667333",
    "2021-10-05 14:34:00.667335",
    "2021-10-05 14:34:00.667338",
    "2021-10-05 14:34:00.667343",
    "2021-10-05 14:34:00.667346",
    "2021-10-05 14:34:00.667348",
    "2021-10-05 14:34:00.667351",
    "2021-10-05 14:34:00.667354",
    "2021-10-05 14:34:00.667356",
    "2021-10-05 14:34:00.667359",
    "2021-10-05 14:34:00.667361",
    "2021-10-05 14:34:00.667364",
    "2021-10-05 14:34:00.667367",
    "2021-10-05 14:34:00.667369",
    "2021-10-05 14:34:00.667372",
    "2021-10-05 14:34:00.667375",
    "2021-10-05 14:34:00.667378",
    "2021-10-05 14:34:00.667381",
    "2021-10-05 14:34:00.667383",
    "2021-10-05 14:34:00.667386",
    "2021-10-05 14:34:00.667389",
    "2021-10-05 14:34:00.667391",
    "2021-10-05 14:34:00.667394",
    "2021-10-05 14:34:00.667397",
    "2021-10-05 14:34:00.667399",
    "2021-10-05 14:34:00.667402",
    "2021-10-05 14:34:00.667405",
    "2021-10-05 14:34:00.667407",
    "2021-10-05 14:34:00.667410",
    "2021-10-05 14:34:00.667412",
    "2021-10-05 14:34:00.667415",
    "2021-10-05 14:34:00.667417",
    "2021-10-05 14:34:00.667420",
    "2021-10-05 14:34:00.667423",
    "2021-10-05 14:34:00.667425",
    "2021-10-05 14:34:00.667428",
    "2021-10-05 14:34:00.667431",
    "2021-10-05 14:34:00.667433",
    "2021-10-05 14:34:00.667436",
    "2021-10-05 14:34:00.667438",
    "2021-10-05 14:34:00.667441",
    "2021-10-05 14:34:00.667444",
    "2021-10-05 14:34:00.667447",
    "2021-10-05 14:34:00.667449",
    "2021-10-05 14:34:00.667452",
    "2021-10-05 14:34:00.667454",
    "2021-10-05 14:34:00.667457",
    "2021-10-05 14:34:00.667461",
    "2021-10-05 14:34:00.667464",
    "2021-10-05 14:34:00.667466",
    "2021-10-05 14:34:00.667469",
    "2021-10-05 14:34:00.667472",
    "2021-10-05 14:34:00.667474",
    "2021-10-05 14:34:00.667477",
    "2021-10-05 14:34:00.667480",
    "2021-10-05 14:34:00.667482",
    "2021-10-05 14:34:00.667485",
    "2021-10-05 14:34:00.667488",
    "2021-10-05 14:34:00.667490"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $newTimestamps[$i]
}

# --- Add a new "metadata" sheet after the "data" sheet ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$meta = $wb.Worksheets.Add($null, $lastSheet)
$meta.Name = "metadata"

# Header row
$meta.Cells.Item(1,2).Value = "data_name"
$meta.Cells.Item(1,3).Value = "data_id"
$meta.Cells.Item(1,4).Value = "data_version"
$meta.Cells.Item(1,5).Value = "data_version_created"
$meta.Cells.Item(1,6).Value = "panel_query_time"
$meta.Cells.Item(1,7).Value = "panel_get_request"

# Data row
$meta.Cells.Item(2,1).Value = 0
$meta.Cells.Item(2,2).Value = "Hereditary Spastic Paraplegia Superpanel"
$meta.Cells.Item(2,3).Value = 3117
$meta.Cells.Item(2,4).NumberFormat = "@"
$meta.Cells.Item(2,4).Value = "2.16"
$meta.Cells.Item(2,5).Value = "2021-10-04T05:19:07.050260Z"
$meta.Cells.Item(2,6).Value = "2021-10-05 14:34:00.663562"
$meta.Cells.Item(2,7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3117/?format=json"

# Match the look of the "data" sheet: bold/bordered header style, and the
# same "index" style used on column A of the data rows.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$ws.Activate()

Write-Output "Edit complete"
